# Add min columns management
#
# Appends three more data rows (9-11) to Sheet1, repeating the same
# values/labels already used in rows 4-6 (AAA/BBB/CCC with their
# Decimal/Integer/Decimal/Data figures). The new rows reuse the existing
# shared strings (no new unique text is introduced) and pick up the
# "plain" cell style (the one already used by columns C/D in row 6 and
# by the blank header row 2) for every numeric column, while the date
# column keeps the existing date-formatted style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (B2:E2) is blank but already carries the "plain" style that the
# new rows need for columns C and D (row 4/5 use a different, italic-ish
# style for those two columns). Row 4's F column already has the date
# style we want to keep for the new F cells.
$plainFormatSource = $ws.Range("C2:D2")

$sourceRows = @(4, 5, 6)
$destRows   = @(9, 10, 11)

for ($i = 0; $i -lt $sourceRows.Length; $i++) {
    $srcRow = $sourceRows[$i]
    $dstRow = $destRows[$i]

    $srcRange = $ws.Range("B" + $srcRow + ":F" + $srcRow)
    $dstRange = $ws.Range("B" + $dstRow + ":F" + $dstRow)

    # Copy values + formulas + styles verbatim from the source row.
    $srcRange.Copy($dstRange)

    # Columns C/D on the source rows use a different style than what the
    # new rows should end up with; overwrite just their formatting (not
    # their values) with the "plain" style from row 2.
    $plainFormatSource.Copy()
    $destFormatTarget = $ws.Range("C" + $dstRow + ":D" + $dstRow)
    $destFormatTarget.PasteSpecial(-4122)
}

Write-Output "Added rows 9-11 (min columns management)"
